# Adding within5 and closest columns to database, making more data persist
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Row 12: FA_within_five_meters
$ws.Range("C12").Value = "integer"
$ws.Range("E12").Value = "FA_within_five_meters"
$ws.Range("F12").Value = "Within 5 meters"

# Row 13: FA_closest_to_focal
$ws.Range("C13").Value = "integer"
$ws.Range("E13").Value = "FA_closest_to_focal"
$ws.Range("F13").Value = "Closest to focal"

# Move the active selection to C14, matching the saved cursor position
[void]$ws.Range("C14").Select()
